$d = $word.ActiveDocument

function Replace-ParagraphXml($findText, $innerXml) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find paragraph containing: $findText"
    }
    $p = $rng.Paragraphs(1)
    $pRng = $p.Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $pRng.InsertXML($xml)
}

# 1) "Z = Multimonster multiplier" -> split runs, proofErr around "Multimonster"
Replace-ParagraphXml "Z = Multimonster multiplier" (
    '<w:r><w:tab/><w:t xml:space="preserve">Z = </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Multimonster</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> multiplier</w:t></w:r>'
)

# 2) "P = Party encounter Xp" -> split runs, proofErr around "Xp"
Replace-ParagraphXml "P = Party encounter Xp" (
    '<w:r><w:tab/><w:t xml:space="preserve">P = Party encounter </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Xp</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>'
)

# 3) "PDeadly = 36x^2 [EN DASH] 99x + 163" -> split runs, proofErr around "PDeadly"
Replace-ParagraphXml "PDeadly = 36x^2" (
    '<w:r><w:tab/></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>PDeadly</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> = 36x^2 ' + [char]0x2013 + ' 99x + 163</w:t></w:r>'
)

# 4) "PNormal = 17x^2 -59x + 92 " -> split runs, proofErr around "PNormal"
Replace-ParagraphXml "PNormal = 17x^2 -59x + 92" (
    '<w:r><w:tab/></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>PNormal</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> = 17x^2 -59x + 92 </w:t></w:r>'
)

# 5) Split the trailing centered bookmark paragraph into an empty centered
#    paragraph followed by a new "Difficulty: 7" paragraph (bookmark moves
#    into the new paragraph).
$last = $d.Paragraphs($d.Paragraphs.Count)
$lastRng = $last.Range
$lastRng.Delete()
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p>' + `
    '<w:p><w:r><w:t>Difficulty:</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> 7</w:t></w:r></w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$lastRng.InsertXML($xml)
